# Apply the "full data scraped for extra batting and bowling fields" update
# to the "ODI Batting Extra" sheet (PlayerPerformance_4063.xlsx / IND team).
#
# The sheet originally only held 20 data rows (rows 2-21). The update
# replaces the data with a fuller 72-row dataset (rows 2-73), re-keyed by
# MATCH_CODE, and refreshes the BATTING_POSITION / NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH values to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Each entry: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$data = @(
    @('3895', '', '', '', '', 'NO'),
    @('3897', '', '', '', '', 'NO'),
    @('3899', '', '', '', '', 'NO'),
    @('4061', 10, '', '', '', 'NO'),
    @('4062', 11, '0', '0', '0.40%', 'NO'),
    @('4063', '', '', '', '', 'NO'),
    @('4065', '', '', '', '', 'NO'),
    @('4067', '', '', '', '', 'NO'),
    @('4069', 11, '', '', '', 'NO'),
    @('4071', 11, '0', '0', '', 'NO'),
    @('4074', '', '', '', '', 'NO'),
    @('4085', '', '', '', '', 'NO'),
    @('4088', 11, '', '', '', 'NO'),
    @('4089', 11, '', '', '', 'NO'),
    @('4096', 11, '', '', '', 'NO'),
    @('4098', '', '', '', '', 'NO'),
    @('4099', '', '', '', '', 'NO'),
    @('4126', 10, '', '', '', 'NO'),
    @('4127', 10, '', '', '', 'NO'),
    @('4128', 9, '', '', '', 'NO'),
    @('4130', '', '', '', '', 'NO'),
    @('4133', 9, '', '', '', 'NO'),
    @('4135', 9, '0', '0', '1.25%', 'NO'),
    @('4171', 9, '', '', '', 'NO'),
    @('4173', '', '', '', '', 'NO'),
    @('4175', 9, '', '', '', 'NO'),
    @('4196', 10, '3', '0', '19.57%', 'NO'),
    @('4197', 11, '0', '0', '', 'NO'),
    @('4199', 10, '0', '0', '', 'NO'),
    @('4201', '', '', '', '', 'NO'),
    @('4205', 10, '', '', '', 'NO'),
    @('4213', 10, '', '', '', 'NO'),
    @('4216', '', '', '', '', 'NO'),
    @('4219', 10, '', '', '', 'NO'),
    @('4236', 10, '', '', '', 'NO'),
    @('4239', 11, '', '', '', 'NO'),
    @('4242', '', '', '', '', 'NO'),
    @('4245', '', '', '', '', 'NO'),
    @('4248', '', '', '', '', 'NO'),
    @('4249', 10, '', '', '', 'NO'),
    @('4268', '', '', '', '', 'NO'),
    @('4310', '', '', '', '', 'NO'),
    @('4316', '', '', '', '', 'NO'),
    @('4324', 11, '', '', '', 'NO'),
    @('4332', 10, '0', '0', '0.71%', 'NO'),
    @('4338', '', '', '', '', 'NO'),
    @('4342', 10, '1', '0', '4.64%', 'NO'),
    @('4345', 10, '', '', '', 'NO'),
    @('4353', 10, '0', '0', '2.05%', 'NO'),
    @('4362', 9, '', '', '', 'NO'),
    @('4406', 10, '', '', '', 'NO'),
    @('4410', '', '', '', '', 'NO'),
    @('4435', 9, '', '', '', 'NO'),
    @('4436', 11, '', '', '', 'NO'),
    @('4480', '', '', '', '', 'NO'),
    @('4482', 9, '0', '0', '3.65%', 'NO'),
    @('4524', '', '', '', '', 'NO'),
    @('4526', 11, '', '', '', 'NO'),
    @('4529', '', '', '', '', ''),
    @('4533', '', '', '', '', ''),
    @('4535', '', '', '', '', ''),
    @('4609', '', '', '', '', ''),
    @('4613', '', '', '', '', ''),
    @('4618', '', '', '', '', ''),
    @('4621', '', '', '', '', ''),
    @('4623', '', '', '', '', ''),
    @('4624', '', '', '', '', ''),
    @('4669', '', '', '', '', ''),
    @('4673', '', '', '', '', ''),
    @('4676', '', '', '', '', ''),
    @('4687', '', '', '', '', ''),
    @('4697', '', '', '', '', '')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $matchCode = $entry[0]
    $battingPosition = $entry[1]
    $num4 = $entry[2]
    $num6 = $entry[3]
    $percentRuns = $entry[4]
    $manOfMatch = $entry[5]

    $ws.Cells.Item($row, 1).Value = $matchCode

    if ($battingPosition -eq '') {
        $ws.Cells.Item($row, 2).ClearContents()
    } else {
        $ws.Cells.Item($row, 2).Value = $battingPosition
    }

    if ($num4 -eq '') {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $num4
    }

    if ($num6 -eq '') {
        $ws.Cells.Item($row, 4).ClearContents()
    } else {
        $ws.Cells.Item($row, 4).Value = $num6
    }

    if ($percentRuns -eq '') {
        $ws.Cells.Item($row, 5).ClearContents()
    } else {
        $ws.Cells.Item($row, 5).Value = $percentRuns
    }

    if ($manOfMatch -eq '') {
        $ws.Cells.Item($row, 6).ClearContents()
    } else {
        $ws.Cells.Item($row, 6).Value = $manOfMatch
    }
}
